$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "55.740.80"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +3.21%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.495.10"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +8.14%  "

$ws.Range("E4").Value = "  +0.00%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "480.02"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +8.21%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "139.27"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +9.05%  "

$ws.Range("E7").Value = "  +0.24%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.513"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +7.91%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.493.51"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +14.06%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0986"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +7.45%  "

$ws.Range("E11").Value = "  +1.53%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.327"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +5.77%  "

$ws.Range("E13").Value = "  +0.37%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.925.62"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +7.64%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "55.764.92"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +3.27%  "

$ws.Range("B16").Value = "Avalanche"
$ws.Range("C16").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "20.49"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +9.33%  "

$ws.Range("B17").Value = "ShibaInu"
$ws.Range("C17").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000137"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +13.79%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.495.36"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +6.74%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.35"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +7.76%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "320.64"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +7.83%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "9.96"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +6.50%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.995"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.31%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.66"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +6.56%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "57.91"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +4.31%  "

$ws.Range("E25").Value = "  +0.88%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.403"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +9.76%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.164"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +5.49%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.608.90"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +8.11%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.39"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +8.04%  "

$ws.Range("E30").Value = "  +8.75%  "

$ws.Range("E31").Value = "  +0.30%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "148.80"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.68%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "18.12"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +7.19%  "

$ws.Range("E34").Value = "  +9.64%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.19"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +11.05%  "

$ws.Range("E36").Value = "  +2.61%  "

$ws.Range("E37").Value = "  +10.48%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.846"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.29%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "34.28"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +3.83%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.999"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.36%  "

$ws.Range("E41").Value = "  +18.29%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0550"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +11.64%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.38"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +7.80%  "

$ws.Range("E44").Value = "  +7.64%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "10.18"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.01%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.969.63"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.40%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0906"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +9.81%  "

$ws.Range("B48").Value = "Bittensor"
$ws.Range("C48").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "251.11"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +33.15%  "

$ws.Range("B49").Value = "VeChain"
$ws.Range("C49").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0222"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +6.81%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "4.50"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +11.37%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "17.41"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +8.53%  "
